$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells D1 and E1 (inherit bold/centered style from row default)
$ws.Range("D1").Value = "ADH_EXP_pred_nocons_bls"
$ws.Range("E1").Value = "ADH_EXP_pred_bls"

# Updated / new numeric values for columns B, C (revised) and D, E (new) across rows 2-51
$ws.Range("B2").Value = 0.01404933931031386
$ws.Range("C2").Value = 0.0132484011039413
$ws.Range("D2").Value = 0.0104079778198634
$ws.Range("E2").Value = 0.009853175537074745
$ws.Range("B3").Value = 0.001015981600812475
$ws.Range("C3").Value = 0.0007606684729950857
$ws.Range("D3").Value = 0.000640974294239542
$ws.Range("E3").Value = 0.0004244143757820513
$ws.Range("B4").Value = 0.009214848000104742
$ws.Range("C4").Value = 0.008920584347595394
$ws.Range("D4").Value = 0.004893050789723919
$ws.Range("E4").Value = 0.004628879292681788
$ws.Range("B5").Value = 0.01139466892848532
$ws.Range("C5").Value = 0.01071537699361355
$ws.Range("D5").Value = 0.007844856435932527
$ws.Range("E5").Value = 0.007249515104500511
$ws.Range("B6").Value = 0.0123636969313909
$ws.Range("C6").Value = 0.01188177379212609
$ws.Range("D6").Value = 0.00711232916864492
$ws.Range("E6").Value = 0.006722986841295486
$ws.Range("B7").Value = 0.00789950518771243
$ws.Range("C7").Value = 0.007544131830237946
$ws.Range("D7").Value = 0.003724351884581838
$ws.Range("E7").Value = 0.00339323930168769
$ws.Range("B8").Value = 0.01055856458746568
$ws.Range("C8").Value = 0.01002330938872306
$ws.Range("D8").Value = 0.006377562480451919
$ws.Range("E8").Value = 0.005970870924020959
$ws.Range("B9").Value = 0.006264578220434048
$ws.Range("C9").Value = 0.005736276733957408
$ws.Range("D9").Value = 0.003954658219529825
$ws.Range("E9").Value = 0.003557966879887878
$ws.Range("B10").Value = 0.005728613849503536
$ws.Range("C10").Value = 0.005424437877269449
$ws.Range("D10").Value = 0.003090601138019741
$ws.Range("E10").Value = 0.002874829040100625
$ws.Range("B11").Value = 0.01483019523778148
$ws.Range("C11").Value = 0.01422437379194086
$ws.Range("D11").Value = 0.009090665985282796
$ws.Range("E11").Value = 0.008630940384025482
$ws.Range("B12").Value = 0.002434731528034411
$ws.Range("C12").Value = 0.002239647671864634
$ws.Range("D12").Value = 0.001250259928587364
$ws.Range("E12").Value = 0.001141551703426879
$ws.Range("B13").Value = 0.009733267834019089
$ws.Range("C13").Value = 0.009370618363536594
$ws.Range("D13").Value = 0.005443221223128375
$ws.Range("E13").Value = 0.005163415377313476
$ws.Range("B14").Value = 0.01009982344900693
$ws.Range("C14").Value = 0.009553908196674834
$ws.Range("D14").Value = 0.00670298478078545
$ws.Range("E14").Value = 0.0062786005022627
$ws.Range("B15").Value = 0.01267642819460724
$ws.Range("C15").Value = 0.01163634989825128
$ws.Range("D15").Value = 0.008425643886056558
$ws.Range("E15").Value = 0.007754813125508958
$ws.Range("B16").Value = 0.01068894146424003
$ws.Range("C16").Value = 0.01005367327089968
$ws.Range("D16").Value = 0.005820276602006415
$ws.Range("E16").Value = 0.005394835714898369
$ws.Range("B17").Value = 0.006518083434129455
$ws.Range("C17").Value = 0.006000123676202413
$ws.Range("D17").Value = 0.004180088908538012
$ws.Range("E17").Value = 0.003742261673949242
$ws.Range("B18").Value = 0.009393672028793879
$ws.Range("C18").Value = 0.008804604176717974
$ws.Range("D18").Value = 0.00744910647825863
$ws.Range("E18").Value = 0.006928756595371975
$ws.Range("B19").Value = 0.004199959668052586
$ws.Range("C19").Value = 0.003705948763951053
$ws.Range("D19").Value = 0.002718475723688833
$ws.Range("E19").Value = 0.002374592097611275
$ws.Range("B20").Value = 0.009960111945414761
$ws.Range("C20").Value = 0.009487525356833027
$ws.Range("D20").Value = 0.006697436111134449
$ws.Range("E20").Value = 0.006314777904538574
$ws.Range("B21").Value = 0.004211235858772567
$ws.Range("C21").Value = 0.00395542515478061
$ws.Range("D21").Value = 0.002801671137003139
$ws.Range("E21").Value = 0.002578702218209118
$ws.Range("B22").Value = 0.01195963632612149
$ws.Range("C22").Value = 0.01146350403643168
$ws.Range("D22").Value = 0.007132278760640352
$ws.Range("E22").Value = 0.006721392517498132
$ws.Range("B23").Value = 0.009163836382168627
$ws.Range("C23").Value = 0.008368765918347209
$ws.Range("D23").Value = 0.0071635853646961
$ws.Range("E23").Value = 0.00655680140505789
$ws.Range("B24").Value = 0.01015894970463629
$ws.Range("C24").Value = 0.00949550165784525
$ws.Range("D24").Value = 0.006256188773145984
$ws.Range("E24").Value = 0.005750121617254366
$ws.Range("B25").Value = 0.01089926505630288
$ws.Range("C25").Value = 0.01014484548719094
$ws.Range("D25").Value = 0.009116557617433038
$ws.Range("E25").Value = 0.008473568829274606
$ws.Range("B26").Value = 0.006932630762336255
$ws.Range("C26").Value = 0.006401497383057692
$ws.Range("D26").Value = 0.005403755124678399
$ws.Range("E26").Value = 0.005002716608703664
$ws.Range("B27").Value = 0.002773783043947493
$ws.Range("C27").Value = 0.002404605466651206
$ws.Range("D27").Value = 0.001546523332050142
$ws.Range("E27").Value = 0.001295201620571976
$ws.Range("B28").Value = 0.006275480363489992
$ws.Range("C28").Value = 0.005853145103836855
$ws.Range("D28").Value = 0.004541878414570482
$ws.Range("E28").Value = 0.004196646120929682
$ws.Range("B29").Value = 0.002582572450243142
$ws.Range("C29").Value = 0.002286212421857174
$ws.Range("D29").Value = 0.001671506636993457
$ws.Range("E29").Value = 0.001436346520318615
$ws.Range("B30").Value = 0.01680707647309597
$ws.Range("C30").Value = 0.01630679471706677
$ws.Range("D30").Value = 0.01008627706100606
$ws.Range("E30").Value = 0.009601279122900234
$ws.Range("B31").Value = 0.00856005957147467
$ws.Range("C31").Value = 0.008042811056598376
$ws.Range("D31").Value = 0.005017335035629103
$ws.Range("E31").Value = 0.004618584495904896
$ws.Range("B32").Value = 0.005050303842318173
$ws.Range("C32").Value = 0.004628996145168004
$ws.Range("D32").Value = 0.003088792425014819
$ws.Range("E32").Value = 0.002813255801954871
$ws.Range("B33").Value = 0.009115261067057914
$ws.Range("C33").Value = 0.008659572061807577
$ws.Range("D33").Value = 0.005452606196576288
$ws.Range("E33").Value = 0.005108976436053196
$ws.Range("B34").Value = 0.02227482857993069
$ws.Range("C34").Value = 0.02138131389035354
$ws.Range("D34").Value = 0.01383363544905441
$ws.Range("E34").Value = 0.01316135730348048
$ws.Range("B35").Value = 0.003959926775825898
$ws.Range("C35").Value = 0.0036476887507078
$ws.Range("D35").Value = 0.002309732946265354
$ws.Range("E35").Value = 0.002067386111202996
$ws.Range("B36").Value = 0.01022205630830632
$ws.Range("C36").Value = 0.009462569411348416
$ws.Range("D36").Value = 0.00724697138305599
$ws.Range("E36").Value = 0.006655491154091646
$ws.Range("B37").Value = 0.008437492914458772
$ws.Range("C37").Value = 0.007782710216582228
$ws.Range("D37").Value = 0.004534899504586422
$ws.Range("E37").Value = 0.004082787612524721
$ws.Range("B38").Value = 0.01122460959563311
$ws.Range("C38").Value = 0.01082938413360711
$ws.Range("D38").Value = 0.006109762857778367
$ws.Range("E38").Value = 0.005747578432740375
$ws.Range("B39").Value = 0.01076204498308192
$ws.Range("C39").Value = 0.01011137891457977
$ws.Range("D39").Value = 0.007340467944595149
$ws.Range("E39").Value = 0.006833545102303642
$ws.Range("B40").Value = 0.01378658440053279
$ws.Range("C40").Value = 0.0129244146149817
$ws.Range("D40").Value = 0.00831953516904162
$ws.Range("E40").Value = 0.007716745967718373
$ws.Range("B41").Value = 0.01965604944678975
$ws.Range("C41").Value = 0.01878862850667255
$ws.Range("D41").Value = 0.01333165117473448
$ws.Range("E41").Value = 0.01275811575018004
$ws.Range("B42").Value = 0.007478081340695028
$ws.Range("C42").Value = 0.007072983089387116
$ws.Range("D42").Value = 0.00511870936894376
$ws.Range("E42").Value = 0.004742293716065738
$ws.Range("B43").Value = 0.01285454350968868
$ws.Range("C43").Value = 0.01204411280438693
$ws.Range("D43").Value = 0.008366785394538253
$ws.Range("E43").Value = 0.007805428799881726
$ws.Range("B44").Value = 0.008850765120638933
$ws.Range("C44").Value = 0.008382396342877966
$ws.Range("D44").Value = 0.005161834184881375
$ws.Range("E44").Value = 0.004775816923017246
$ws.Range("B45").Value = 0.008068843455869751
$ws.Range("C45").Value = 0.007481202230758182
$ws.Range("D45").Value = 0.004415765717841744
$ws.Range("E45").Value = 0.003973386880619191
$ws.Range("B46").Value = 0.011280632616462
$ws.Range("C46").Value = 0.01076392706601916
$ws.Range("D46").Value = 0.003878039959123733
$ws.Range("E46").Value = 0.003388304795514888
$ws.Range("B47").Value = 0.007767335858419836
$ws.Range("C47").Value = 0.007346530104131937
$ws.Range("D47").Value = 0.00488952137096842
$ws.Range("E47").Value = 0.004538953336002602
$ws.Range("B48").Value = 0.006055920202743124
$ws.Range("C48").Value = 0.005646955813134055
$ws.Range("D48").Value = 0.003786344020771993
$ws.Range("E48").Value = 0.003433140915458698
$ws.Range("B49").Value = 0.00583852661062775
$ws.Range("C49").Value = 0.005138487772201776
$ws.Range("D49").Value = 0.003272073411936396
$ws.Range("E49").Value = 0.002775377249096812
$ws.Range("B50").Value = 0.01275775208597906
$ws.Range("C50").Value = 0.01196612010467505
$ws.Range("D50").Value = 0.008345517900613001
$ws.Range("E50").Value = 0.007725213939243207
$ws.Range("B51").Value = 0.001661664059564683
$ws.Range("C51").Value = 0.001103771846363116
$ws.Range("D51").Value = 0.001363815388569738
$ws.Range("E51").Value = 0.0009536728367573426
